{"js": "// Values taken from the target revision: the banner date plus every\n// arithmetic-problem cell of the practice table, in row-major order\n// (20 rows x 5 columns), exactly matching the source document's structure.\nconst dateNew = \"2024-10-04 Friday\";\nconst newValues = [[\"90-8=\", \"49+26=\", \"73+20=\", \"92-86=\", \"93-50=\"], [\"20+63=\", \"48-40=\", \"48+48=\", \"38+5=\", \"39-7=\"], [\"56-4=\", \"65-46=\", \"42-37=\", \"80-30=\", \"4+68=\"], [\"38-5=\", \"35+59=\", \"12+45=\", \"26+45=\", \"64+4=\"], [\"20+63=\", \"38+2=\", \"36+32=\", \"53+3=\", \"68-54=\"], [\"18+14=\", \"98-39=\", \"40-29=\", \"7+64=\", \"28+70=\"], [\"0+48=\", \"50+6=\", \"3+16=\", \"28+71=\", \"96-89=\"], [\"87-47=\", \"7+48=\", \"29+5=\", \"9+86=\", \"70-41=\"], [\"95-48=\", \"91-51=\", \"77-47=\", \"86+3=\", \"85-79=\"], [\"72-12=\", \"98-74=\", \"57-21=\", \"48-40=\", \"96-16=\"], [\"71-38=\", \"41+34=\", \"43-40=\", \"29+36=\", \"53-47=\"], [\"75-7=\", \"39-15=\", \"79-2=\", \"73+11=\", \"37+9=\"], [\"76-41=\", \"27+24=\", \"3+11=\", \"52-2=\", \"30+1=\"], [\"35+16=\", \"41-29=\", \"88-50=\", \"82-12=\", \"77-22=\"], [\"60-25=\", \"78-71=\", \"33+36=\", \"90-28=\", \"44+24=\"], [\"91-26=\", \"4+5=\", \"79-44=\", \"53+32=\", \"53+9=\"], [\"85-16=\", \"33+55=\", \"27+4=\", \"93-9=\", \"93-72=\"], [\"60+4=\", \"25+23=\", \"7+1=\", \"73-24=\", \"13+80=\"], [\"79-41=\", \"62-27=\", \"48+27=\", \"80-17=\", \"76-51=\"], [\"71+22=\", \"53+32=\", \"89-58=\", \"48-5=\", \"38+14=\"]];\n\nconst body = context.document.body;\n\n// Update the date line (first paragraph in the document body). Using the\n// paragraph's own range keeps the existing run formatting (font/size).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.getRange(\"Whole\").insertText(dateNew, Word.InsertLocation.replace);\n\n// Update every cell of the table, addressed by (row, column) so that\n// duplicate problem text (e.g. two different cells that both happened to\n// read \"92+1=\") is replaced with the correct, position-specific value\n// instead of being matched by content.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let r = 0; r < newValues.length; r++) {\n  const rowValues = newValues[r];\n  for (let c = 0; c < rowValues.length; c++) {\n    const cell = table.getCell(r, c);\n    const cellParagraph = cell.body.paragraphs.getFirst();\n    cellParagraph.getRange(\"Whole\").insertText(rowValues[c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date line (first paragraph in the document body).\n$d.Paragraphs.Item(1).Range.Text = \"2024-10-04 Friday\"\n\n# Update each arithmetic-problem cell in the table, in row-major order,\n# matching the structure of the source document exactly.\n$t = $d.Tables.Item(1)\n\n$values = @(\n    @(\"90-8=\", \"49+26=\", \"73+20=\", \"92-86=\", \"93-50=\"),\n    @(\"20+63=\", \"48-40=\", \"48+48=\", \"38+5=\", \"39-7=\"),\n    @(\"56-4=\", \"65-46=\", \"42-37=\", \"80-30=\", \"4+68=\"),\n    @(\"38-5=\", \"35+59=\", \"12+45=\", \"26+45=\", \"64+4=\"),\n    @(\"20+63=\", \"38+2=\", \"36+32=\", \"53+3=\", \"68-54=\"),\n    @(\"18+14=\", \"98-39=\", \"40-29=\", \"7+64=\", \"28+70=\"),\n    @(\"0+48=\", \"50+6=\", \"3+16=\", \"28+71=\", \"96-89=\"),\n    @(\"87-47=\", \"7+48=\", \"29+5=\", \"9+86=\", \"70-41=\"),\n    @(\"95-48=\", \"91-51=\", \"77-47=\", \"86+3=\", \"85-79=\"),\n    @(\"72-12=\", \"98-74=\", \"57-21=\", \"48-40=\", \"96-16=\"),\n    @(\"71-38=\", \"41+34=\", \"43-40=\", \"29+36=\", \"53-47=\"),\n    @(\"75-7=\", \"39-15=\", \"79-2=\", \"73+11=\", \"37+9=\"),\n    @(\"76-41=\", \"27+24=\", \"3+11=\", \"52-2=\", \"30+1=\"),\n    @(\"35+16=\", \"41-29=\", \"88-50=\", \"82-12=\", \"77-22=\"),\n    @(\"60-25=\", \"78-71=\", \"33+36=\", \"90-28=\", \"44+24=\"),\n    @(\"91-26=\", \"4+5=\", \"79-44=\", \"53+32=\", \"53+9=\"),\n    @(\"85-16=\", \"33+55=\", \"27+4=\", \"93-9=\", \"93-72=\"),\n    @(\"60+4=\", \"25+23=\", \"7+1=\", \"73-24=\", \"13+80=\"),\n    @(\"79-41=\", \"62-27=\", \"48+27=\", \"80-17=\", \"76-51=\"),\n    @(\"71+22=\", \"53+32=\", \"89-58=\", \"48-5=\", \"38+14=\")\n)\n\nfor ($r = 0; $r -lt $values.Count; $r++) {\n    $rowValues = $values[$r]\n    for ($c = 0; $c -lt $rowValues.Count; $c++) {\n        $t.Cell($r + 1, $c + 1).Range.Text = $rowValues[$c]\n    }\n}\n\n"}
